# Applies the "Add files via upload" revision to Final_Equipment_List.xlsx
#
# Content changes:
#  1. Equipment-ReactorSection: P-601 A/B pump "Power (motor)" label gets a
#     space before "kW".
#  2. Equipment-ReactorSection: E-601 Dowtherm A Cooler duty/heat-transfer
#     area values are updated (4376 -> 4400 MJ/h, 5.2 -> 5.4 m^2).
#  3. Heat Exchangers: the various "Stainless Steel" / "Monel" material
#     call-outs are consolidated into a single "Monel/Carbon Steel" label.
#  4. Selection / active-sheet bookkeeping: Equipment-ReactorSection's
#     selection moves to M10, the Heat Exchangers tab becomes the active
#     sheet (with selection Q19), and Compressors is no longer the active
#     tab.

$wb = $excel.ActiveWorkbook

# --- 1 & 2: Equipment-ReactorSection sheet -------------------------------
$reactor = $wb.Worksheets.Item("Equipment-ReactorSection")

$reactor.Range("A21").Value = "Power (motor) = 2.5 kW"
$reactor.Range("F7").Value = "Heat Transfer Area = 5.4 m^2"
$reactor.Range("F6").Value = "Duty = 4400 MJ/h"

# --- 3: Heat Exchangers sheet --------------------------------------------
$heat = $wb.Worksheets.Item("Heat Exchangers")

$heat.Range("A3").Value = "Monel/Carbon Steel"
$heat.Range("A9").Value = "Monel/Carbon Steel"
$heat.Range("A15").Value = "Monel/Carbon Steel"
$heat.Range("A21").Value = "Monel/Carbon Steel"
$heat.Range("A27").Value = "Monel/Carbon Steel"
$heat.Range("A33").Value = "Monel/Carbon Steel"

# A3 / A9 previously carried stray direct formatting with no visible
# effect (plain / italic font on top of the default style); clear it so
# the cells fall back to the default, unformatted style.
$heat.Range("A3").Font.Bold = $false
$heat.Range("A9").Font.Italic = $false

# --- 4: selection / active sheet bookkeeping -----------------------------
$reactor.Range("M10").Select() | Out-Null

$heat.Activate()
$heat.Range("Q19").Select() | Out-Null
